$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 46216.46301711571
$ws.Range("C2").Value = 21711.32539286149
$ws.Range("D2").Value = 58211.43633129547
